# Updated cryptos list on Mon Jun 10 20:52:27 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D = Price, Column E = Volume(1h). Both columns hold text (not
# numeric) values in this sheet, so every Price cell we touch gets its
# NumberFormat forced to "@" (text) before the assignment - otherwise
# Excel auto-coerces plain-looking decimals (e.g. "626.16") into numbers.

function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
}

# Row 2 - Bitcoin
Set-TextCell "D2" "69.652.12"
$ws.Range("E2").Value = "  +0.02%  "

# Row 3 - Ethereum
Set-TextCell "D3" "3.673.35"
$ws.Range("E3").Value = "  -0.68%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.08%  "

# Row 5 - BNB
Set-TextCell "D5" "626.16"
$ws.Range("E5").Value = "  -6.93%  "

# Row 6 - Solana
Set-TextCell "D6" "160.23"
$ws.Range("E6").Value = "  -0.72%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.01%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -0.05%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -2.01%  "

# Row 10 - Toncoin
Set-TextCell "D10" "7.16"
$ws.Range("E10").Value = "  +0.99%  "

# Row 11 - Cardano
Set-TextCell "D11" "0.440"
$ws.Range("E11").Value = "  -1.08%  "

# Row 12 - ShibaInu
Set-TextCell "D12" "0.0000229"

# Row 13 - WrappedliquidstakedEther2.0
Set-TextCell "D13" "4.292.16"
$ws.Range("E13").Value = "  -0.69%  "

# Row 14 - Avalanche
Set-TextCell "D14" "32.55"
$ws.Range("E14").Value = "  -0.72%  "

# Row 15 - WrappedEther
Set-TextCell "D15" "3.678.70"
$ws.Range("E15").Value = "  -0.64%  "

# Row 16 - WrappedBTC
Set-TextCell "D16" "69.674.24"
$ws.Range("E16").Value = "  +0.04%  "

# Row 17 - TRON
$ws.Range("E17").Value = "  +0.76%  "

# Row 18 - Polkadot
Set-TextCell "D18" "6.53"
$ws.Range("E18").Value = "  +0.36%  "

# Row 19 - Chainlink
Set-TextCell "D19" "15.90"
$ws.Range("E19").Value = "  -1.68%  "

# Row 20 - Uniswap
Set-TextCell "D20" "10.35"
$ws.Range("E20").Value = "  +5.53%  "

# Row 21 - BitcoinCash
Set-TextCell "D21" "471.10"
$ws.Range("E21").Value = "  -0.63%  "

# Row 22 - Polygon
Set-TextCell "D22" "0.651"

# Row 23 - Litecoin
Set-TextCell "D23" "79.77"
$ws.Range("E23").Value = "  -0.87%  "

# Row 24 - WrappedeETH
Set-TextCell "D24" "3.819.46"
$ws.Range("E24").Value = "  -0.74%  "

# Row 25 - Dai
$ws.Range("E25").Value = "  +0.12%  "

# Row 26 - PEPE
$ws.Range("E26").Value = "  -1.70%  "

# Row 27 - InternetComputer(DFINITY)
Set-TextCell "D27" "11.07"

# Row 28 - RenderToken
Set-TextCell "D28" "8.72"
$ws.Range("E28").Value = "  -4.51%  "

# Row 29 - PancakeSwap
Set-TextCell "D29" "2.59"
$ws.Range("E29").Value = "  -3.73%  "

# Rows 31, 32, 34 - reorder (ImmutableX, EthereumClassic, Binance-PegBSC-USD
# were re-ranked) plus value updates. Row 33 (Kaspa) stays in place.

# Row 31 -> Binance-PegBSC-USD
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextCell "D31" "1.01"
$ws.Range("E31").Value = "  +0.73%  "

# Row 32 -> ImmutableX
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell "D32" "1.99"
$ws.Range("E32").Value = "  -1.47%  "

# Row 33 - Kaspa (name/link/price unchanged, only Volume(1h) changes)
$ws.Range("E33").Value = "  -0.39%  "

# Row 34 -> EthereumClassic
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell "D34" "26.65"
$ws.Range("E34").Value = "  -0.84%  "

# Row 35 - NEARProtocol
$ws.Range("E35").Value = "  -2.87%  "

# Row 36 - RenzoRestakedETH
Set-TextCell "D36" "3.675.78"
$ws.Range("E36").Value = "  -0.32%  "

# Row 37 - Aptos
Set-TextCell "D37" "8.32"
$ws.Range("E37").Value = "  -2.43%  "

# Row 38 - USDe
$ws.Range("E38").Value = "  -0.05%  "

# Row 39 - Monero
Set-TextCell "D39" "178.45"
$ws.Range("E39").Value = "  +3.30%  "

# Row 40 - Filecoin
Set-TextCell "D40" "5.83"
$ws.Range("E40").Value = "  -4.98%  "

# Row 41 - FirstDigitalUSD
Set-TextCell "D41" "0.999"

# Row 42 - Stacks
Set-TextCell "D42" "2.19"
$ws.Range("E42").Value = "  -1.39%  "

# Row 43 - Hedera
$ws.Range("E43").Value = "  -1.27%  "

# Row 44 - Mantle
$ws.Range("E44").Value = "  -1.48%  "

# Row 45 - OKB
Set-TextCell "D45" "46.75"
$ws.Range("E45").Value = "  -0.55%  "

# Row 46 - InjectiveProtocol
$ws.Range("E46").Value = "  +3.82%  "

# Row 47 - dogwifhat
Set-TextCell "D47" "2.73"
$ws.Range("E47").Value = "  -1.52%  "

# Row 48 - Cosmos
Set-TextCell "D48" "7.89"
$ws.Range("E48").Value = "  -0.04%  "

# Row 49 - FLOKI
$ws.Range("E49").Value = "  -4.80%  "

# Row 50 - SuiNetwork
$ws.Range("E50").Value = "  -5.63%  "

# Row 51 - ONDO
$ws.Range("E51").Value = "  -5.33%  "
